$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6: convert the "weight" column (C) from a plain inline-string
# number (e.g. "2", "10", "20") to a real numeric percentage value
# (e.g. 0.02, 0.10, 0.20), clear the now-unused "Your input" column (D)
# entirely, and strip the explicit formatting from the result column (E).

$ws.Range("C2").Value = 0.02
$ws.Range("D2").Clear()
$ws.Range("E2").ClearFormats()

$ws.Range("C3").Value = 0.1
$ws.Range("D3").Clear()
$ws.Range("E3").ClearFormats()

$ws.Range("C4").Value = 0.2
$ws.Range("D4").Clear()
$ws.Range("E4").ClearFormats()

$ws.Range("C5").Value = 0.2
$ws.Range("D5").Clear()
$ws.Range("E5").ClearFormats()

$ws.Range("C6").Value = 0.48
$ws.Range("D6").Clear()
$ws.Range("E6").ClearFormats()

# New row 7: a "Total" label plus a formula that sums the per-row totals.
$ws.Range("D7").Value = "Total"
$ws.Range("E7").Formula = "=E2+E3+E4+E5+E6"
$ws.Range("E7").ClearFormats()
